$wb = $excel.ActiveWorkbook
Write-Output "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
